$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.663177523997797
$ws.Range("C2").Value = 0.663177523997797
$ws.Range("D2").Value = 4.2919259948383
$ws.Range("E2").Value = 0.0227109217402772
$ws.Range("F2").Value = 0.0003

$ws.Range("B3").Value = 5.79639801105541
$ws.Range("C3").Value = 1.9321326703518
$ws.Range("D3").Value = 12.5042995778409
$ws.Range("E3").Value = 0.198501210974417

$ws.Range("B4").Value = 1.10879878102872
$ws.Range("C4").Value = 0.369599593676241
$ws.Range("D4").Value = 2.39195999016698
$ws.Range("E4").Value = 0.0379714954600026
$ws.Range("F4").Value = 0.0014

$ws.Range("B5").Value = 21.6324450774202
$ws.Range("C5").Value = 0.154517464838716
$ws.Range("E5").Value = 0.740816371825303

$ws.Range("B6").Value = 29.2008193935022
